$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 2816.625
$ws.Range("I34").Value = 2816.625
$ws.Range("K34").Value = 2816.625
$ws.Range("M34").Value = -2613.625
$ws.Range("H36").Value = 2816.625
$ws.Range("I36").Value = 2816.625
$ws.Range("K36").Value = 2816.625
$ws.Range("M36").Value = -2101.625
$ws.Range("H39").Value = 498.35715
$ws.Range("I39").Value = 186.77777
$ws.Range("J39").Value = 1059.2
$ws.Range("K39").Value = 560.33331
$ws.Range("L39").Value = 3177.6
$ws.Range("M39").Value = -264.33331
$ws.Range("N39").Value = -3769.6
$ws.Range("H43").Value = 5350.591
$ws.Range("I43").Value = 2021.5714
$ws.Range("J43").Value = 6904.1333
$ws.Range("K43").Value = 2021.5714
$ws.Range("L43").Value = 6904.1333
$ws.Range("M43").Value = -1952.5714
$ws.Range("N43").Value = -7042.1333
$ws.Range("H54").Value = 204494.8
$ws.Range("I54").Value = 204494.8
$ws.Range("K54").Value = 204494.8
$ws.Range("M54").Value = -204008.8
$ws.Range("H133").Value = 69508
$ws.Range("J133").Value = 69508
$ws.Range("L133").Value = 69508
$ws.Range("N133").Value = -79628
$ws.Range("H134").Value = 69999
$ws.Range("J134").Value = 69999
$ws.Range("L134").Value = 69999
$ws.Range("H138").Value = 2791.3193
$ws.Range("I138").Value = 1057
$ws.Range("K138").Value = 3171
$ws.Range("M138").Value = 1969
$ws.Range("H141").Value = 2351.182
$ws.Range("I141").Value = 2006.4445
$ws.Range("K141").Value = 6019.333500000001
$ws.Range("M141").Value = -839.3335000000006
$ws.Range("N134").Value = -80139

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2657.2786
$ws.Range("I32").Value = 2235.9868
$ws.Range("J32").Value = 13330
$ws.Range("K32").Value = 2235.9868
$ws.Range("L32").Value = 13330
$ws.Range("M32").Value = -1948.9868
$ws.Range("N32").Value = -13904
$ws.Range("H74").Value = 11114276
$ws.Range("I74").Value = 15153464
$ws.Range("K74").Value = 15153464
$ws.Range("M74").Value = -15152590
$ws.Range("H77").Value = 11114276
$ws.Range("I77").Value = 15153464
$ws.Range("K77").Value = 75767320
$ws.Range("M77").Value = -75762952
$ws.Range("H135").Value = 56047.547
$ws.Range("J135").Value = 56047.547
$ws.Range("L135").Value = 56047.547
$ws.Range("N135").Value = -66187.54699999999
$ws.Range("H137").Value = 69999
$ws.Range("J137").Value = 69999
$ws.Range("L137").Value = 69999
$ws.Range("N137").Value = -80199
$ws.Range("H139").Value = 58956.445
$ws.Range("J139").Value = 58956.445
$ws.Range("L139").Value = 58956.445
$ws.Range("N139").Value = -69236.44500000001
$ws.Range("H141").Value = 57848
$ws.Range("J141").Value = 57848
$ws.Range("L141").Value = 57848
$ws.Range("N141").Value = -68208

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 5000
$ws.Range("I33").Value = 5000
$ws.Range("K33").Value = 5000
$ws.Range("H97").Value = 20175.25
$ws.Range("I97").Value = 5356.25
$ws.Range("K97").Value = 5356.25
$ws.Range("M97").Value = -4365.25
$ws.Range("H132").Value = 69999
$ws.Range("J132").Value = 69999
$ws.Range("L132").Value = 69999
$ws.Range("H135").Value = 46879.25
$ws.Range("J135").Value = 46879.25
$ws.Range("L135").Value = 46879.25
$ws.Range("N135").Value = -57019.25
$ws.Range("H137").Value = 69999
$ws.Range("J137").Value = 69999
$ws.Range("L137").Value = 69999
$ws.Range("H138").Value = 59391.223
$ws.Range("J138").Value = 59391.223
$ws.Range("L138").Value = 59391.223
$ws.Range("N138").Value = -69671.223
$ws.Range("H140").Value = 68992.8
$ws.Range("J140").Value = 68992.8
$ws.Range("L140").Value = 68992.8
$ws.Range("N140").Value = -79352.8
$ws.Range("M33").Value = -4664
$ws.Range("N132").Value = -80119
$ws.Range("N137").Value = -80199

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 100
$ws.Range("I21").Value = 100
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 100
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 135
$ws.Range("H22").Value = 1365.4546
$ws.Range("I22").Value = 420
$ws.Range("K22").Value = 420
$ws.Range("M22").Value = -70
$ws.Range("H31").Value = 33046.51
$ws.Range("I31").Value = 3664.88
$ws.Range("J31").Value = 94258.25
$ws.Range("K31").Value = 3664.88
$ws.Range("L31").Value = 94258.25
$ws.Range("M31").Value = -3369.88
$ws.Range("N31").Value = -94848.25
$ws.Range("H32").Value = 4941.8
$ws.Range("I32").Value = 1569.6666
$ws.Range("J32").Value = 10000
$ws.Range("K32").Value = 1569.6666
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = -1253.6666
$ws.Range("N32").Value = -10632
$ws.Range("H34").Value = 33046.51
$ws.Range("I34").Value = 3664.88
$ws.Range("J34").Value = 94258.25
$ws.Range("K34").Value = 3664.88
$ws.Range("L34").Value = 94258.25
$ws.Range("M34").Value = -3462.88
$ws.Range("N34").Value = -94662.25
$ws.Range("H35").Value = 1131.1765
$ws.Range("I35").Value = 718.9231
$ws.Range("J35").Value = 2471
$ws.Range("K35").Value = 718.9231
$ws.Range("L35").Value = 2471
$ws.Range("M35").Value = -424.9231
$ws.Range("H38").Value = 2500
$ws.Range("J38").Value = 2500
$ws.Range("L38").Value = 2500
$ws.Range("N38").Value = -3254
$ws.Range("H46").Value = 2500
$ws.Range("J46").Value = 2500
$ws.Range("L46").Value = 2500
$ws.Range("N46").Value = -2922
$ws.Range("H62").Value = 10198.125
$ws.Range("I62").Value = 3695
$ws.Range("K62").Value = 3695
$ws.Range("M62").Value = -3071
$ws.Range("H65").Value = 10198.125
$ws.Range("I65").Value = 3695
$ws.Range("K65").Value = 18475
$ws.Range("M65").Value = -15355
$ws.Range("H93").Value = 10140.8
$ws.Range("I93").Value = 10140.8
$ws.Range("K93").Value = 10140.8
$ws.Range("M93").Value = -8268.799999999999
$ws.Range("H99").Value = 2959.7
$ws.Range("I99").Value = 2259.4
$ws.Range("K99").Value = 2259.4
$ws.Range("M99").Value = -761.4000000000001
$ws.Range("H103").Value = 37497.5
$ws.Range("I103").Value = 37497.5
$ws.Range("K103").Value = 37497.5
$ws.Range("M103").Value = -36325.5
$ws.Range("H107").Value = 1559.0869
$ws.Range("I107").Value = 1456.3125
$ws.Range("J107").Value = 1794
$ws.Range("K107").Value = 1456.3125
$ws.Range("L107").Value = 1794
$ws.Range("M107").Value = 463.6875
$ws.Range("N107").Value = -5634
$ws.Range("H126").Value = 2959.7
$ws.Range("I126").Value = 2259.4
$ws.Range("K126").Value = 6778.200000000001
$ws.Range("M126").Value = -4308.200000000001
$ws.Range("H132").Value = 2696.9666
$ws.Range("I132").Value = 1756.3636
$ws.Range("J132").Value = 5283.625
$ws.Range("K132").Value = 5269.0908
$ws.Range("L132").Value = 15850.875
$ws.Range("M132").Value = -2739.0908
$ws.Range("N132").Value = -20910.875
$ws.Range("H135").Value = 69992.8
$ws.Range("J135").Value = 69992.8
$ws.Range("L135").Value = 69992.8
$ws.Range("N135").Value = -80132.8
$ws.Range("H138").Value = 64188.8
$ws.Range("J138").Value = 64188.8
$ws.Range("L138").Value = 64188.8
$ws.Range("N138").Value = -74468.8
$ws.Range("H140").Value = 68776
$ws.Range("J140").Value = 68776
$ws.Range("L140").Value = 68776
$ws.Range("N140").Value = -79136
$ws.Range("N21").ClearContents()
$ws.Range("N35").Value = -3059

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5348214
$ws.Range("I4").Value = 3461839.5
$ws.Range("J4").Value = 9806916
$ws.Range("K4").Value = 10385518.5
$ws.Range("L4").Value = 29420748
$ws.Range("M4").Value = -10385406.5
$ws.Range("N4").Value = -29420972
$ws.Range("H5").Value = 1964099.9
$ws.Range("I5").Value = 2231.6667
$ws.Range("J5").Value = 3034209.8
$ws.Range("K5").Value = 6695.000100000001
$ws.Range("L5").Value = 9102629.399999999
$ws.Range("M5").Value = -6583.000100000001
$ws.Range("N5").Value = -9102853.399999999
$ws.Range("H94").Value = 898.3333
$ws.Range("I94").Value = 898.3333
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 2694.9999
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -2018.9999
$ws.Range("H129").Value = 7577236
$ws.Range("J129").Value = 41670170
$ws.Range("L129").Value = 125010510
$ws.Range("N129").Value = -125020510
$ws.Range("H130").Value = 1197.5
$ws.Range("I130").Value = 1197.5
$ws.Range("K130").Value = 3592.5
$ws.Range("M130").Value = 1427.5
$ws.Range("H135").Value = 1964099.9
$ws.Range("I135").Value = 2231.6667
$ws.Range("J135").Value = 3034209.8
$ws.Range("K135").Value = 20085.0003
$ws.Range("L135").Value = 27307888.2
$ws.Range("M135").Value = -17550.0003
$ws.Range("N135").Value = -27312958.2
$ws.Range("H136").Value = 2839.3333
$ws.Range("I136").Value = 2214.8572
$ws.Range("J136").Value = 5025
$ws.Range("K136").Value = 6644.571599999999
$ws.Range("L136").Value = 15075
$ws.Range("M136").Value = -1544.571599999999
$ws.Range("N136").Value = -25275
$ws.Range("H137").Value = 102959.9
$ws.Range("J137").Value = 203599.8
$ws.Range("L137").Value = 610799.3999999999
$ws.Range("N137").Value = -620999.3999999999
$ws.Range("H140").Value = 2070.1667
$ws.Range("I140").Value = 1165.4
$ws.Range("K140").Value = 3496.2
$ws.Range("M140").Value = 1683.8
$ws.Range("H141").Value = 8773.8125
$ws.Range("I141").Value = 4542.857
$ws.Range("J141").Value = 12064.556
$ws.Range("K141").Value = 13628.571
$ws.Range("L141").Value = 36193.66800000001
$ws.Range("M141").Value = -8448.571
$ws.Range("N141").Value = -46553.66800000001
$ws.Range("N94").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 606.3333
$ws.Range("J2").Value = 1227.1428
$ws.Range("L2").Value = 1227.1428
$ws.Range("N2").Value = -1453.1428
$ws.Range("H80").Value = 4995.5557
$ws.Range("I80").Value = 2868.75
$ws.Range("J80").Value = 6697
$ws.Range("K80").Value = 2868.75
$ws.Range("L80").Value = 6697
$ws.Range("M80").Value = -1870.75
$ws.Range("N80").Value = -8693
$ws.Range("H83").Value = 4995.5557
$ws.Range("I83").Value = 2868.75
$ws.Range("J83").Value = 6697
$ws.Range("K83").Value = 14343.75
$ws.Range("L83").Value = 33485
$ws.Range("M83").Value = -9351.75
$ws.Range("N83").Value = -43469
$ws.Range("H93").Value = 29584.416
$ws.Range("J93").Value = 29855.7
$ws.Range("L93").Value = 29855.7
$ws.Range("N93").Value = -33599.7
$ws.Range("H97").Value = 862.93335
$ws.Range("I97").Value = 926.4211
$ws.Range("K97").Value = 926.4211
$ws.Range("M97").Value = -430.4211
$ws.Range("H102").Value = 3750.963
$ws.Range("I102").Value = 2740.65
$ws.Range("K102").Value = 2740.65
$ws.Range("M102").Value = -1118.65
$ws.Range("H122").Value = 5568.5454
$ws.Range("I122").Value = 4434.55
$ws.Range("K122").Value = 13303.65
$ws.Range("M122").Value = -10853.65
$ws.Range("H132").Value = 5109.4062
$ws.Range("I132").Value = 3847.3914
$ws.Range("K132").Value = 11542.1742
$ws.Range("M132").Value = -9012.174199999999
$ws.Range("H133").Value = 69999
$ws.Range("J133").Value = 69999
$ws.Range("L133").Value = 69999
$ws.Range("H135").Value = 69999
$ws.Range("J135").Value = 69999
$ws.Range("L135").Value = 69999
$ws.Range("H138").Value = 64987.5
$ws.Range("J138").Value = 64987.5
$ws.Range("L138").Value = 64987.5
$ws.Range("N138").Value = -75267.5
$ws.Range("H140").Value = 71565.75
$ws.Range("J140").Value = 71565.75
$ws.Range("L140").Value = 71565.75
$ws.Range("N140").Value = -81925.75
$ws.Range("N133").Value = -80119
$ws.Range("N135").Value = -80139

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 23855
$ws.Range("J4").Value = 23855
$ws.Range("L4").Value = 23855
$ws.Range("N4").Value = -24081
$ws.Range("H9").Value = 1624.75
$ws.Range("I9").Value = 166.33333
$ws.Range("K9").Value = 166.33333
$ws.Range("M9").Value = 57.66667000000001
$ws.Range("H11").Value = 2950
$ws.Range("I11").Value = 900
$ws.Range("K11").Value = 900
$ws.Range("H28").Value = 23855
$ws.Range("J28").Value = 23855
$ws.Range("L28").Value = 23855
$ws.Range("N28").Value = -24319
$ws.Range("H37").Value = 23855
$ws.Range("J37").Value = 23855
$ws.Range("L37").Value = 23855
$ws.Range("N37").Value = -24069
$ws.Range("H82").Value = 8589.134
$ws.Range("I82").Value = 4964.357
$ws.Range("K82").Value = 4964.357
$ws.Range("M82").Value = -4603.357
$ws.Range("H85").Value = 8589.134
$ws.Range("I85").Value = 4964.357
$ws.Range("K85").Value = 4964.357
$ws.Range("M85").Value = -3716.357
$ws.Range("H122").Value = 194663.28
$ws.Range("I122").Value = 268833.8
$ws.Range("J122").Value = 9237
$ws.Range("K122").Value = 806501.3999999999
$ws.Range("L122").Value = 27711
$ws.Range("M122").Value = -804051.3999999999
$ws.Range("N122").Value = -32611
$ws.Range("H132").Value = 7920.8213
$ws.Range("I132").Value = 6873.5557
$ws.Range("J132").Value = 9805.9
$ws.Range("K132").Value = 20620.6671
$ws.Range("L132").Value = 29417.7
$ws.Range("M132").Value = -18090.6671
$ws.Range("N132").Value = -34477.7
$ws.Range("H133").Value = 49998.668
$ws.Range("I133").Value = 49996
$ws.Range("J133").Value = 50000
$ws.Range("K133").Value = 49996
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -55060
$ws.Range("H134").Value = 51847
$ws.Range("J134").Value = 51847
$ws.Range("L134").Value = 51847
$ws.Range("N134").Value = -61987
$ws.Range("H137").Value = 65212
$ws.Range("I137").Value = 50890
$ws.Range("J137").Value = 69986
$ws.Range("K137").Value = 50890
$ws.Range("L137").Value = 69986
$ws.Range("M137").Value = -45790
$ws.Range("N137").Value = -80186
$ws.Range("H139").Value = 69987
$ws.Range("J139").Value = 69987
$ws.Range("L139").Value = 69987
$ws.Range("N139").Value = -80267
$ws.Range("H141").Value = 69982
$ws.Range("J141").Value = 69982
$ws.Range("L141").Value = 69982
$ws.Range("N141").Value = -80342
$ws.Range("M11").Value = -760
$ws.Range("M133").Value = -47466

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 679000
$ws.Range("I29").Value = 48750
$ws.Range("J29").Value = 3200000
$ws.Range("K29").Value = 48750
$ws.Range("L29").Value = 3200000
$ws.Range("M29").Value = -48460
$ws.Range("N29").Value = -3200580
$ws.Range("H31").Value = 20000
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("H74").Value = 19251.666
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 19251.666
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 19251.666
$ws.Range("N74").Value = -21123.666
$ws.Range("H77").Value = 19251.666
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 19251.666
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 57754.99800000001
$ws.Range("N77").Value = -67114.99800000001
$ws.Range("H81").Value = 4166.6665
$ws.Range("I81").Value = 4166.6665
$ws.Range("K81").Value = 8333.333000000001
$ws.Range("M81").Value = -7272.333000000001
$ws.Range("H84").Value = 4166.6665
$ws.Range("I84").Value = 4166.6665
$ws.Range("K84").Value = 41666.665
$ws.Range("M84").Value = -36362.665
$ws.Range("H92").Value = 31249.875
$ws.Range("J92").Value = 31249.875
$ws.Range("L92").Value = 31249.875
$ws.Range("N92").Value = -36241.875
$ws.Range("H122").Value = 2459.3713
$ws.Range("I122").Value = 1828.2963
$ws.Range("K122").Value = 5484.8889
$ws.Range("M122").Value = -3034.8889
$ws.Range("H132").Value = 4590.1665
$ws.Range("I132").Value = 2503.3333
$ws.Range("J132").Value = 10850.667
$ws.Range("K132").Value = 7509.999899999999
$ws.Range("L132").Value = 32552.001
$ws.Range("M132").Value = -4979.999899999999
$ws.Range("N132").Value = -37612.001
$ws.Range("H135").Value = 56939.5
$ws.Range("J135").Value = 56939.5
$ws.Range("L135").Value = 56939.5
$ws.Range("N135").Value = -67079.5
$ws.Range("H136").Value = 1861.2413
$ws.Range("J136").Value = 7168
$ws.Range("L136").Value = 21504
$ws.Range("N136").Value = -26604
$ws.Range("H137").Value = 67447.8
$ws.Range("J137").Value = 67447.8
$ws.Range("L137").Value = 67447.8
$ws.Range("N137").Value = -77647.8
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("H139").Value = 69783.09
$ws.Range("J139").Value = 70071
$ws.Range("L139").Value = 70071
$ws.Range("N139").Value = -80351
$ws.Range("H141").Value = 100517.37
$ws.Range("J141").Value = 100517.37
$ws.Range("L141").Value = 100517.37
$ws.Range("N31").ClearContents()
$ws.Range("M57").ClearContents()
$ws.Range("M74").ClearContents()
$ws.Range("M77").ClearContents()
$ws.Range("N138").ClearContents()
